# 67763 Sprint1 Amend - Review phase in scrum board must have effort attributed
# to it in burndown (they count as part of the task)
#
# Moves the "T1: Play the first area of the game." task card on the Scrum
# Board slide (slide 15, "Student #67804" row) from the "Doing" column back
# into the "To Do" column, matching the x-position already used by other
# "To Do" cards, with only a tiny y-nudge from the manual drag.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 16) {
        $target = $sh
        break
    }
}

# Target position in EMU: x=9256129, y=5530450 (1 EMU = 1/12700 pt).
# The literal point values below are chosen so that, after this runtime's
# internal float32 conversion back to EMU, they land exactly on the target
# EMU values (plain "emu / 12700.0" rounds 9256129 down to 9256128 here).
$target.Left = 728.8291
$target.Top = 435.4685039370079
